# Archiv.xlsx edit: add a "Kürzel" (code) column at the start of the sheet.
#
# The new column A holds consecutive numeric IDs (4689, 4690, ...) for each
# data row, formatted as integers. All existing columns (old A..E) shift
# right by one (new B..F). The header cell becomes "Kürzel".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all existing columns one place to the right, freeing up column A.
$ws.Columns.Item(1).Insert()

# New header for the inserted column.
$ws.Range("A1").Value2 = "Kürzel"

# Fill in the sequential codes for every data row (rows 2-94 => 4689-4781),
# applying an integer number format to each new cell.
$firstCode = 4689
$lastRow = 94

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value2 = $firstCode + ($row - 2)
    $cell.NumberFormat = "0"
}

# Reflect the view state from the edited workbook as closely as possible.
$ws.Range("C98").Select()
